$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update maxCount (column C) values: most cards 3 -> 2, a few 3 -> 1 ---
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 2

# --- Update card effect text (column D) ---
# Row 10 "Genocide": shorten effect text
$ws.Range("D10").Value = "消耗3时间：获得1技能点。"

# Row 8 "Midas Touch": shorten effect text and tweak wording
$ws.Range("D8").Value = "消耗3时间：将主牌堆第1张怪物牌放在房间区任意非空列最前方，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆顶的1张遗物牌。"

# --- Update the selected/active cell in the sheet view ---
$ws.Range("D9").Select() | Out-Null
